$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.177.59"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +3.12%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.736.67"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("E4").Value = "  +0.27%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "240.38"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.9980"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4789"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -1.95%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2593"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.44%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06148"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.735.78"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.97%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "16.11"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +3.85%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.06925"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.6020"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.431"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "76.76"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.42%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.9983"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "27.137.98"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.41%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.9977"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.27%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000007074"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.11%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "11.41"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.90%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.946.57"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.65%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.417"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.83%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "8.382"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("E24").Value = "  +1.74%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "141.88"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +4.12%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "15.23"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.43%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.819"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +5.34%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "106.77"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.35%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.381"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.70%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "3.944"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +1.80%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.07922"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.666"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.96%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.04758"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +7.05%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.590"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.56%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.012"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +2.12%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.6166"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.40%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.9210"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.85%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.533"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +7.74%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.017"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +1.36%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.9981"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.11%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.685"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +5.33%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.01488"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.05%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "98.83"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.60%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.3820"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.61%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "6.844"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("E46").Value = "  +0.40%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.05350"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.830"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +2.42%  "

$ws.Range("E49").Value = "  -1.64%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.244"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +3.66%  "

$ws.Range("E51").Value = "  +0.01%  "
